$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Cornflower blue (FF6495ED) expressed as a COM BGR integer, matching the
# workbook's existing custom "HyperLink" cell style color.
$linkColor = 15570276

$commit = "77b07f7431702b20fa613b003e585b02ef206db6"
$guid = "1f5eeef4-23f6-43cc-b531-a6f094206bcb"

$zhXlfName = "$guid.$commit.zh-cn.xlf"
$deXlfName = "$guid.$commit.de-de.xlf"

$zhXlfUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/Loc/zh-cn/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/Loc/de-de/$deXlfName"

# --- Overview sheet: the shared "Handoff failed" label is now reported as
#     "Not yet handed off" everywhere it is used, including here. ---
$wsOverview.Range("B2").Value2 = "Not yet handed off"
$wsOverview.Range("C2").Value2 = "Not yet handed off"

# --- zh-cn sheet updates ---
$wsZh.Range("B2").Value2 = "Not yet handed off"

# Re-create the hyperlinks so the new handoff-file link lands between the
# existing source-file link (A2) and the config link (A3), matching the
# relationship id ordering of the generated report.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3cd67f8cd2b56ef92c6bdb711b23c876d2dc3982/e2e/$guid.md", "", "", "$guid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3cd67f8cd2b56ef92c6bdb711b23c876d2dc3982/.localization-config", "", "", ".localization-config")

# Hyperlinks.Add() re-derives the font/style for touched cells; restore the
# workbook's custom hyperlink look (underline + cornflower blue) on all three.
$wsZh.Range("A2").Font.Underline = $true
$wsZh.Range("A2").Font.Color = $linkColor
$wsZh.Range("A3").Font.Underline = $true
$wsZh.Range("A3").Font.Color = $linkColor
$wsZh.Range("C2").Font.Underline = $true
$wsZh.Range("C2").Font.Color = $linkColor

$wsZh.Range("D2").Value2 = "2016-01-11 03:29:05"
$wsZh.Range("H2").Value2 = "Include"

# --- de-de sheet updates (mirrors zh-cn) ---
$wsDe.Range("B2").Value2 = "Not yet handed off"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3cd67f8cd2b56ef92c6bdb711b23c876d2dc3982/e2e/$guid.md", "", "", "$guid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3cd67f8cd2b56ef92c6bdb711b23c876d2dc3982/.localization-config", "", "", ".localization-config")

$wsDe.Range("A2").Font.Underline = $true
$wsDe.Range("A2").Font.Color = $linkColor
$wsDe.Range("A3").Font.Underline = $true
$wsDe.Range("A3").Font.Color = $linkColor
$wsDe.Range("C2").Font.Underline = $true
$wsDe.Range("C2").Font.Color = $linkColor

$wsDe.Range("D2").Value2 = "2016-01-11 03:29:21"
$wsDe.Range("H2").Value2 = "Include"

Write-Output "Report generated for handoff"
